$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the cells we are updating so that numeric-looking
# strings (e.g. "1.001", "213.33") are preserved verbatim as text instead of
# being auto-converted to numbers by Excel, matching the original inlineStr cells.
$cellUpdates = @{
    "D2" = "26.068.22"
    "E2" = "  -0.20%  "
    "D3" = "1.635.83"
    "E3" = "  -1.90%  "
    "E4" = "  -0.11%  "
    "D5" = "213.33"
    "E5" = "  +1.33%  "
    "D6" = "0.5238"
    "E6" = "  -0.29%  "
    "E7" = "  -0.12%  "
    "E8" = "  -1.12%  "
    "D9" = "0.06307"
    "E9" = "  +0.26%  "
    "D10" = "20.70"
    "E10" = "  -2.30%  "
    "D11" = "0.07650"
    "E11" = "  +1.33%  "
    "D12" = "1.633.99"
    "E12" = "  -2.28%  "
    "D13" = "4.413"
    "E13" = "  -0.54%  "
    "D14" = "1.859.72"
    "E14" = "  -1.92%  "
    "D15" = "0.5492"
    "E15" = "  -0.87%  "
    "D16" = "0.0₅8168"
    "E16" = "  +3.09%  "
    "D17" = "64.95"
    "E17" = "  -2.69%  "
    "D18" = "26.059.15"
    "E18" = "  -0.37%  "
    "E19" = "  -0.04%  "
    "D20" = "4.684"
    "E20" = "  -1.13%  "
    "D21" = "188.13"
    "E21" = "  +0.97%  "
    "D22" = "10.14"
    "E22" = "  -1.86%  "
    "D23" = "6.144"
    "E23" = "  -0.44%  "
    "E24" = "  -0.02%  "
    "D25" = "145.67"
    "E25" = "  -2.78%  "
    "E26" = "  -3.00%  "
    "D27" = "7.391"
    "E27" = "  -1.38%  "
    "D28" = "15.78"
    "E28" = "  -0.92%  "
    "D29" = "1.401"
    "E29" = "  +3.52%  "
    "D30" = "0.05962"
    "E30" = "  -4.56%  "
    "D31" = "1.253"
    "E31" = "  -1.97%  "
    "D32" = "3.433"
    "E32" = "  -2.28%  "
    "D33" = "3.403"
    "E33" = "  -0.29%  "
    "D34" = "1.636"
    "E34" = "  +0.55%  "
    "D35" = "0.9849"
    "E35" = "  -1.20%  "
    "D36" = "2.397"
    "E36" = "  -0.74%  "
    "E37" = "  +1.06%  "
    "D38" = "0.5718"
    "E38" = "  -5.40%  "
    "D39" = "0.01617"
    "E39" = "  +0.23%  "
    "D40" = "0.8530"
    "E40" = "  -2.03%  "
    "B41" = "PaxDollar"
    "C41" = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
    "D41" = "1.001"
    "E41" = "  -0.22%  "
    "B42" = "FraxShare"
    "C42" = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
    "D42" = "5.735"
    "E42" = "  -6.40%  "
    "D43" = "1.034.36"
    "E43" = "  -6.37%  "
    "D44" = "100.54"
    "E44" = "  +0.52%  "
    "D45" = "1.786.10"
    "E45" = "  -1.85%  "
    "D46" = "0.0₈105"
    "E46" = "  -4.65%  "
    "D47" = "55.78"
    "E47" = "  +0.75%  "
    "E48" = "  -0.37%  "
    "D49" = "8.031"
    "E49" = "  -0.31%  "
    "D50" = "0.05166"
    "E50" = "  -1.36%  "
    "D51" = "0.4222"
    "E51" = "  -0.57%  "
}

foreach ($addr in $cellUpdates.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cellUpdates[$addr]
    $rng.Style = "Normal"
}
